$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet2: reset its data so it mirrors Sheet1 ("Clear" button resets the
# lower-case demo data back to the same values/format used on Sheet1),
# then re-apply the few cells that were hand edited afterwards.
# ---------------------------------------------------------------------------

# Row 1 (merged header)
$ws2.Range("A1").Value = $ws1.Range("A1").Value()

# Row 2
$ws2.Range("A2").Value = $ws1.Range("A2").Value()
$ws2.Range("B2").Value = $ws1.Range("B2").Value()
$ws2.Range("C2").Value = $ws1.Range("C2").Value()
$ws2.Range("D2").Value = $ws1.Range("D2").Value()

# Row 3 (merged halves)
$ws2.Range("A3").Value = $ws1.Range("A3").Value()
$ws2.Range("C3").Value = $ws1.Range("C3").Value()

# Row 4
$ws2.Range("A4").Value = $ws1.Range("A4").Value()
$ws2.Range("B4").Value = $ws1.Range("B4").Value()
$ws2.Range("C4").Value = $ws1.Range("C4").Value()
$ws2.Range("D4").Value = "Datac4"

# Row 5
$ws2.Range("A5").Value = $ws1.Range("A5").Value()
$ws2.Range("B5").Value = $ws1.Range("B5").Value()
$ws2.Range("C5").Value = $ws1.Range("C5").Value()
$ws2.Range("D5").Value = $ws1.Range("D5").Value()

# Row 7: checkmark / mark row, rearranged compared to Sheet1's layout
$ws2.Range("B7").ClearContents()
$ws2.Range("A7").Value = $ws1.Range("A7").Value()
$ws2.Range("C7").Value = $ws1.Range("C7").Value()
$ws2.Range("D7").Value = $ws1.Range("B7").Value()

# ---------------------------------------------------------------------------
# Formatting: copy the exact formats from Sheet1's matching cells so that
# the same underlying style records get reused.
# ---------------------------------------------------------------------------
$ws1.Range("A1:D1").Copy()
$ws2.Range("A1:D1").PasteSpecial(-4122)

$ws1.Range("A3:D3").Copy()
$ws2.Range("A3:D3").PasteSpecial(-4122)

$ws1.Range("A7:C7").Copy()
$ws2.Range("A7").PasteSpecial(-4122)
$ws2.Range("C7").PasteSpecial(-4122)
$ws2.Range("D7").PasteSpecial(-4122)

$ws1.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet2 page setup (clear button also set a printable page size once data
# was reset) and selections left on each sheet after the edit.
# ---------------------------------------------------------------------------
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

[void]$ws1.Range("A1:D7").Select()
[void]$ws1.Range("D7").Activate()

[void]$ws2.Range("D4").Select()
